$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-08-19 Tuesday" "2025-08-20 Wednesday"

Replace-Text "155×9=" "481×4="
Replace-Text "875×6=" "656×4="
Replace-Text "191×4=" "124×6="
Replace-Text "271×5=" "497×4="
Replace-Text "513×9=" "349×9="

Replace-Text "569×3=" "851×5="
Replace-Text "627×6=" "868×5="
Replace-Text "564×5=" "723×2="
Replace-Text "224×8=" "449×6="
Replace-Text "658×4=" "964×3="

Replace-Text "892×4=" "867×4="
Replace-Text "914×5=" "418×7="
Replace-Text "585×9=" "146×4="
Replace-Text "722×6=" "411×3="
Replace-Text "347×9=" "973×9="

Replace-Text "441×9=" "594×3="
Replace-Text "275×3=" "735×2="
Replace-Text "915×7=" "518×2="
Replace-Text "737×5=" "714×3="
Replace-Text "684×5=" "290×7="

Replace-Text "919×9=" "137×9="
Replace-Text "691×6=" "436×5="
Replace-Text "507×5=" "296×5="
Replace-Text "284×5=" "658×3="
Replace-Text "175×4=" "440×3="
